# Apply the edit described by the diff:
#  1. Shared string "MODEL_CONDITION" -> "MODELCONDITION"
#  2. Remove the original column A (the redundant row-index column:
#     values 1 / 8), shifting every other column one place to the left
#     (B->A, C->B, D->C, E->D, F->E). This also drops the header-style
#     (s="1") that used to sit on A2/A3, and shrinks the sheet's used
#     range/dimension from A1:F3 to A1:E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the header text for the MODEL_CONDITION column (now column D
#    once the index column is removed, but the text can be fixed first).
$xlWhole = 1
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION", $xlWhole)

# 2) Delete the leftmost column (old column A, holding the redundant
#    row-index values 1 and 8) so every other column shifts one place
#    to the left: B->A, C->B, D->C, E->D, F->E.
$ws.Columns.Item(1).Delete()
